$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Effect Date* for the data row moves forward to 2025-09-29 (stored as text, matching
# the existing quote-prefixed text style already applied to I2).
$ws.Range("I2").Value = "'2025-09-29"

# Refreshed price columns for the row
$ws.Range("J2").Value = 9879.0
$ws.Range("K2").Value = 9879.0
$ws.Range("L2").Value = 9880.0
$ws.Range("M2").Value = 9880.0
